# Meal.xlsx maintenance edit:
# Insert a header row at the top of Sheet1 (mealName / special / type),
# shifting all existing meal data down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing rows down and create a blank row 1.
[void]$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "mealName"
$ws.Range("B1").Value = "special"
$ws.Range("C1").Value = "type"

# Match the saved selection state (cell B1 active).
[void]$ws.Range("B1").Select()
